$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add new row 8: A8 = "a" (same string as A2/A3), B8 = 3.5
$ws.Range("A8").Value = "a"
$ws.Range("B8").Value = 3.5

# Update the active selection to B9, matching the post-edit state
$ws.Range("B9").Select()
